{"js": "// Update the worksheet date and the twenty-five two-digit multiplication\n// problems to the new day's values. Every old value is unique in the\n// document, so an exact, case-sensitive text search-and-replace is safe.\nconst replacements = [\n  [\"2026-01-17 Saturday\", \"2026-01-18 Sunday\"],\n  [\"21\u00d776=\", \"54\u00d744=\"],\n  [\"54\u00d791=\", \"33\u00d765=\"],\n  [\"38\u00d787=\", \"46\u00d718=\"],\n  [\"61\u00d758=\", \"70\u00d755=\"],\n  [\"43\u00d714=\", \"92\u00d713=\"],\n  [\"36\u00d748=\", \"76\u00d718=\"],\n  [\"25\u00d798=\", \"83\u00d723=\"],\n  [\"43\u00d785=\", \"82\u00d749=\"],\n  [\"24\u00d722=\", \"71\u00d758=\"],\n  [\"25\u00d716=\", \"89\u00d714=\"],\n  [\"74\u00d790=\", \"89\u00d735=\"],\n  [\"27\u00d731=\", \"73\u00d777=\"],\n  [\"13\u00d753=\", \"15\u00d750=\"],\n  [\"91\u00d731=\", \"48\u00d794=\"],\n  [\"65\u00d775=\", \"58\u00d750=\"],\n  [\"84\u00d793=\", \"29\u00d723=\"],\n  [\"68\u00d753=\", \"64\u00d780=\"],\n  [\"92\u00d789=\", \"23\u00d764=\"],\n  [\"31\u00d728=\", \"92\u00d738=\"],\n  [\"17\u00d798=\", \"86\u00d765=\"],\n  [\"32\u00d774=\", \"18\u00d779=\"],\n  [\"89\u00d791=\", \"99\u00d787=\"],\n  [\"55\u00d753=\", \"43\u00d775=\"],\n  [\"25\u00d792=\", \"74\u00d740=\"],\n  [\"12\u00d722=\", \"38\u00d774=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twenty-five two-digit multiplication\n# problems to the new day's values. Every old value is unique in the\n# document, so an exact, case-sensitive find-and-replace-all is safe.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"2026-01-17 Saturday\" = \"2026-01-18 Sunday\"\n    \"21\u00d776=\" = \"54\u00d744=\"\n    \"54\u00d791=\" = \"33\u00d765=\"\n    \"38\u00d787=\" = \"46\u00d718=\"\n    \"61\u00d758=\" = \"70\u00d755=\"\n    \"43\u00d714=\" = \"92\u00d713=\"\n    \"36\u00d748=\" = \"76\u00d718=\"\n    \"25\u00d798=\" = \"83\u00d723=\"\n    \"43\u00d785=\" = \"82\u00d749=\"\n    \"24\u00d722=\" = \"71\u00d758=\"\n    \"25\u00d716=\" = \"89\u00d714=\"\n    \"74\u00d790=\" = \"89\u00d735=\"\n    \"27\u00d731=\" = \"73\u00d777=\"\n    \"13\u00d753=\" = \"15\u00d750=\"\n    \"91\u00d731=\" = \"48\u00d794=\"\n    \"65\u00d775=\" = \"58\u00d750=\"\n    \"84\u00d793=\" = \"29\u00d723=\"\n    \"68\u00d753=\" = \"64\u00d780=\"\n    \"92\u00d789=\" = \"23\u00d764=\"\n    \"31\u00d728=\" = \"92\u00d738=\"\n    \"17\u00d798=\" = \"86\u00d765=\"\n    \"32\u00d774=\" = \"18\u00d779=\"\n    \"89\u00d791=\" = \"99\u00d787=\"\n    \"55\u00d753=\" = \"43\u00d775=\"\n    \"25\u00d792=\" = \"74\u00d740=\"\n    \"12\u00d722=\" = \"38\u00d774=\"\n}\n\nforeach ($oldText in $replacements.Keys) {\n    $newText = $replacements[$oldText]\n    $findRange = $d.Content\n    $findRange.Find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
